# Fixed data errors for citations and language toggle #124
#
# The Spanish citation path in "pages_with_translation"!A9 was truncated/
# wrong ("...examenes-de"); correct it to the full path ending in
# "...examenes-de-deteccion". Excel stores this as a (new) shared string,
# so simply writing the corrected text to the cell creates it.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # pages_with_translation
$ws2 = $wb.Worksheets.Item(2)   # pages_without_translation

# --- Correct the truncated Spanish citation URL in A9 -----------------
$ws1.Range("A9").Value = "espanol/noticias/comunicados-de-prensa/2019/aprendizaje-profundo-cancer-cuello-uterino-examenes-de-deteccion"

# --- Column A on sheet 1 needs to be much wider to show the corrected,
#     much longer URL without truncation. -------------------------------
$ws1.Columns.Item(1).ColumnWidth = 99.6667   # -> stored width 100.5

# --- Row 7 no longer needs its custom row height; restore it to the
#     sheet's standard height. ------------------------------------------
$ws1.Rows.Item(7).AutoFit()

# --- Leave the cursor / selection on the cell that was corrected. ------
$ws1.Activate()
$ws1.Range("A9").Select()
